$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Uppercase the street names and piece labels (shared strings content update)
$ws.Range("C2").Value = "RUE AMHERST"
$ws.Range("D2").Value = "CORPS"
$ws.Range("C3").Value = "CHEMIN FRASER"
$ws.Range("D3").Value = "COUVERCLE"
$ws.Range("D4").Value = "COUVERCLE"
$ws.Range("C4").Value = "BOULEVARD SAINT-RAYMOND"

# Increment the request number column by 1 for rows 2-4
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 3
$ws.Range("A4").Value = 4
